$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.765.84"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.701.23"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").Value = "'316.09"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").Value = "'0.3938"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.4060"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "'1.507"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'1.009"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").Value = "'52.40"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "'0.08783"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").Value = "'7.554"
$ws.Range("E13").Value = "  +5.21%  "
$ws.Range("D14").Value = "'24.56"
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("D15").Value = "'0.00001365"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "'7.999"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "1.703.49"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'99.40"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "'0.07097"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").Value = "'19.87"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "'7.388"
$ws.Range("E21").Value = "  +5.29%  "
$ws.Range("D22").Value = "'1.012"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'14.33"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "24.762.45"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Value = "'3.039"
$ws.Range("E25").Value = "  -5.58%  "
$ws.Range("D26").Value = "'2.350"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'22.74"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'164.95"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'8.639"
$ws.Range("E29").Value = "  +15.80%  "
$ws.Range("D30").Value = "'138.08"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "'5.226"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.891.03"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'7.636"
$ws.Range("E33").Value = "  +6.11%  "
$ws.Range("D34").Value = "'0.08872"
$ws.Range("E34").Value = "  +3.27%  "
$ws.Range("D35").Value = "'1.045"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").Value = "'1.989"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("D37").Value = "'0.2740"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +6.19%  "
$ws.Range("E39").Value = "  -5.72%  "
$ws.Range("D40").Value = "'14.35"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'0.09123"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "'0.7837"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("D43").Value = "'1.470"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'16.55"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "'0.7224"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'2.594"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "'4.247"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D49").Value = "'1.328"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "'139.90"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'91.87"
$ws.Range("E51").Value = "  +2.56%  "
